$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 761.8333
$ws.Range("I2").Value = 810.8
$ws.Range("K2").Value = 810.8
$ws.Range("M2").Value = -697.8

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 448.45456
$ws.Range("I18").Value = 463.4
$ws.Range("J18").Value = 299
$ws.Range("K18").Value = 463.4
$ws.Range("L18").Value = 299
$ws.Range("M18").Value = -179.4
$ws.Range("N18").Value = -867

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 1955.7778
$ws.Range("J70").Value = 1640
$ws.Range("L70").Value = 4920
$ws.Range("N70").Value = -5460

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H73").Value = 1955.7778
$ws.Range("J73").Value = 1640
$ws.Range("L73").Value = 4920
$ws.Range("N73").Value = -6792

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 209437.67
$ws.Range("J75").Value = 209437.67
$ws.Range("L75").Value = 209437.67
$ws.Range("N75").Value = -211309.67

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H78").Value = 209437.67
$ws.Range("J78").Value = 209437.67
$ws.Range("L78").Value = 628313.01
$ws.Range("N78").Value = -637673.01

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 50003804
$ws.Range("J112").Value = 50003804
$ws.Range("L112").Value = 150011412
$ws.Range("N112").Value = -150013628

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12422.087
$ws.Range("I32").Value = 11783.35
$ws.Range("J32").Value = 16680.334
$ws.Range("K32").Value = 11783.35
$ws.Range("L32").Value = 16680.334
$ws.Range("M32").Value = -11496.35
$ws.Range("N32").Value = -17254.334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5944.222
$ws.Range("J45").Value = 6199.6
$ws.Range("L45").Value = 6199.6
$ws.Range("N45").Value = -6953.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1499.4286
$ws.Range("I97").Value = 1315.25
$ws.Range("J97").Value = 2088.8
$ws.Range("K97").Value = 1315.25
$ws.Range("L97").Value = 2088.8
$ws.Range("M97").Value = -819.25
$ws.Range("N97").Value = -3080.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H74").Value = 29664.25
$ws.Range("I74").Value = 29400
$ws.Range("J74").Value = 29752.334
$ws.Range("K74").Value = 29400
$ws.Range("L74").Value = 29752.334
$ws.Range("M74").Value = -28464
$ws.Range("N74").Value = -31624.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H77").Value = 29664.25
$ws.Range("I77").Value = 29400
$ws.Range("J77").Value = 29752.334
$ws.Range("K77").Value = 88200
$ws.Range("L77").Value = 89257.00199999999
$ws.Range("M77").Value = -83520
$ws.Range("N77").Value = -98617.00199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 8884.846
$ws.Range("I105").Value = 13262.875
$ws.Range("K105").Value = 13262.875
$ws.Range("M105").Value = -11515.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2612
$ws.Range("I132").Value = 839.6
$ws.Range("J132").Value = 5566
$ws.Range("K132").Value = 2518.8
$ws.Range("L132").Value = 16698
$ws.Range("M132").Value = 11.19999999999982
$ws.Range("N132").Value = -21758

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 277499.38
$ws.Range("J141").Value = 277499.38
$ws.Range("L141").Value = 277499.38
$ws.Range("N141").Value = -287859.38

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 777.8
$ws.Range("I12").Value = 344
$ws.Range("K12").Value = 1032
$ws.Range("M12").Value = -859

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 2118.8125
$ws.Range("I68").Value = 1258.4
$ws.Range("J68").Value = 2509.9092
$ws.Range("K68").Value = 3775.2
$ws.Range("L68").Value = 7529.7276
$ws.Range("M68").Value = -2964.2
$ws.Range("N68").Value = -9151.7276

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 2118.8125
$ws.Range("I71").Value = 1258.4
$ws.Range("J71").Value = 2509.9092
$ws.Range("K71").Value = 11325.6
$ws.Range("L71").Value = 22589.1828
$ws.Range("M71").Value = -7269.6
$ws.Range("N71").Value = -30701.1828

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 9000
$ws.Range("J106").Value = 9000
$ws.Range("L106").Value = 27000
$ws.Range("N106").Value = -28892

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 31129.578
$ws.Range("J131").Value = 4451.645
$ws.Range("L131").Value = 13354.935
$ws.Range("N131").Value = -23434.935

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 7203.222
$ws.Range("J134").Value = 14999.667
$ws.Range("L134").Value = 44999.001
$ws.Range("N134").Value = -55139.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 3480.4546
$ws.Range("I139").Value = 3531.2856
$ws.Range("J139").Value = 3391.5
$ws.Range("K139").Value = 10593.8568
$ws.Range("L139").Value = 10174.5
$ws.Range("M139").Value = -5453.856800000001
$ws.Range("N139").Value = -20454.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 12158.895
$ws.Range("I141").Value = 10632.182
$ws.Range("J141").Value = 14258.125
$ws.Range("K141").Value = 31896.546
$ws.Range("L141").Value = 42774.375
$ws.Range("M141").Value = -26716.546
$ws.Range("N141").Value = -53134.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 776.6667
$ws.Range("I31").Value = 776.6667
$ws.Range("K31").Value = 776.6667
$ws.Range("M31").Value = -484.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H37").Value = 776.6667
$ws.Range("I37").Value = 776.6667
$ws.Range("K37").Value = 776.6667
$ws.Range("M37").Value = -499.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3047.8948
$ws.Range("I80").Value = 3002.7273
$ws.Range("J80").Value = 3110
$ws.Range("K80").Value = 3002.7273
$ws.Range("L80").Value = 3110
$ws.Range("M80").Value = -2004.7273
$ws.Range("N80").Value = -5106

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3047.8948
$ws.Range("I83").Value = 3002.7273
$ws.Range("J83").Value = 3110
$ws.Range("K83").Value = 15013.6365
$ws.Range("L83").Value = 15550
$ws.Range("M83").Value = -10021.6365
$ws.Range("N83").Value = -25534

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2662.4119
$ws.Range("J97").Value = 3781.8
$ws.Range("L97").Value = 3781.8
$ws.Range("N97").Value = -4773.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2596.84
$ws.Range("I102").Value = 1456.5
$ws.Range("J102").Value = 4624.1113
$ws.Range("K102").Value = 1456.5
$ws.Range("L102").Value = 4624.1113
$ws.Range("M102").Value = 165.5
$ws.Range("N102").Value = -7868.1113

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 24924.857
$ws.Range("J109").Value = 22395
$ws.Range("L109").Value = 22395
$ws.Range("N109").Value = -24475

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1932.6666
$ws.Range("I132").Value = 1771.2858
$ws.Range("K132").Value = 5313.857400000001
$ws.Range("M132").Value = -2783.857400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1012.5
$ws.Range("J22").Value = 2300
$ws.Range("L22").Value = 2300
$ws.Range("N22").Value = -2890

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1012.5
$ws.Range("J27").Value = 2300
$ws.Range("L27").Value = 2300
$ws.Range("N27").Value = -2514

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H54").Value = 35998
$ws.Range("I54").Value = 39995
$ws.Range("J54").Value = 33999.5
$ws.Range("K54").Value = 39995
$ws.Range("L54").Value = 33999.5
$ws.Range("M54").Value = -39351
$ws.Range("N54").Value = -35287.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2716.7144
$ws.Range("I82").Value = 1683.909
$ws.Range("K82").Value = 1683.909
$ws.Range("M82").Value = -1322.909

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2716.7144
$ws.Range("I85").Value = 1683.909
$ws.Range("K85").Value = 1683.909
$ws.Range("M85").Value = -435.9090000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3797.2273
$ws.Range("I100").Value = 3323.7144
$ws.Range("J100").Value = 4625.875
$ws.Range("K100").Value = 3323.7144
$ws.Range("L100").Value = 4625.875
$ws.Range("M100").Value = -2782.7144
$ws.Range("N100").Value = -5707.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6331.8887
$ws.Range("I136").Value = 5666
$ws.Range("K136").Value = 16998
$ws.Range("M136").Value = -14448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 27000
$ws.Range("I51").Value = 24500
$ws.Range("K51").Value = 24500
$ws.Range("M51").Value = -23990

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 24747
$ws.Range("I52").Value = 8999
$ws.Range("K52").Value = 8999
$ws.Range("M52").Value = -8773

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 80000
$ws.Range("J70").Value = 80000
$ws.Range("L70").Value = 80000
$ws.Range("N70").Value = -80630

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 80000
$ws.Range("J73").Value = 80000
$ws.Range("L73").Value = 80000
$ws.Range("N73").Value = -82184

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1442
$ws.Range("I81").Value = 1367.3334
$ws.Range("J81").Value = 1591.3334
$ws.Range("K81").Value = 2734.6668
$ws.Range("L81").Value = 3182.6668
$ws.Range("M81").Value = -1673.6668
$ws.Range("N81").Value = -5304.6668

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1442
$ws.Range("I84").Value = 1367.3334
$ws.Range("J84").Value = 1591.3334
$ws.Range("K84").Value = 13673.334
$ws.Range("L84").Value = 15913.334
$ws.Range("M84").Value = -8369.333999999999
$ws.Range("N84").Value = -26521.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1166.3334
$ws.Range("I113").Value = 750
$ws.Range("K113").Value = 2250
$ws.Range("M113").Value = -80

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 19666.666
$ws.Range("I126").Value = 19666.666
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 58999.99800000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -56529.99800000001
$ws.Range("N126").ClearContents()
